$wb = $excel.ActiveWorkbook

# --- Sheet "Stundenerfassung" (sheet1): add two new rows of logged work ---
$ws1 = $wb.Worksheets.Item("Stundenerfassung")

# Row 125 - copy the date style from the row above first so the new date
# cell keeps the existing "date" cell style (instead of Excel inventing a
# brand-new number format / style record for it).
$ws1.Range("A124").Copy($ws1.Range("A125"))
$ws1.Cells.Item(125, 1).Value = 42975
$ws1.Cells.Item(125, 2).Value = "ETIC2"
$ws1.Cells.Item(125, 3).Value = "Codierung nach MVVM"
$ws1.Cells.Item(125, 4).Value = 2

# Row 126
$ws1.Range("A124").Copy($ws1.Range("A126"))
$ws1.Cells.Item(126, 1).Value = 42975
$ws1.Cells.Item(126, 2).Value = "Schriftliche Arbeit"
$ws1.Cells.Item(126, 3).Value = "Fertigstellung der Arbeit"
$ws1.Cells.Item(126, 4).Value = 8

$ws1.Application.CutCopyMode = $false

# Update the view: this sheet becomes the active one, scrolled further down,
# with a single-cell selection on E124.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 114
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("E124").Select() | Out-Null

# --- Sheet "Wochen" (sheet3): no longer the active tab, scrolled to A7 ---
$ws3 = $wb.Worksheets.Item("Wochen")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("E15").Select() | Out-Null

# Re-activate "Stundenerfassung" so it ends up as the selected / active tab
$ws1.Activate()
